# chore: update Sheets via scheduled runner
# Refreshes cached market-price / profit figures (currentAveragePrice*,
# LevePriceNQ/HQ, LeveProfitNQ/HQ) across the per-job leve-profit tables.
$wb = $excel.ActiveWorkbook

# ALC row 6
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 49.5
$ws.Range("I6").Value = 49.5
$ws.Range("K6").Value = 148.5
$ws.Range("M6").Value = -36.5

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2205.8635
$ws.Range("J17").Value = 2244.238
$ws.Range("L17").Value = 6732.714
$ws.Range("N17").Value = -7068.714

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 3525.3333
$ws.Range("I112").Value = 3600
$ws.Range("K112").Value = 10800
$ws.Range("M112").Value = -9692

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 902.5714
$ws.Range("I132").Value = 902.5714
$ws.Range("K132").Value = 2707.7142
$ws.Range("M132").Value = -177.7142000000003

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1855.7273
$ws.Range("I137").Value = 1802.25
$ws.Range("J137").Value = 1998.3334
$ws.Range("K137").Value = 5406.75
$ws.Range("L137").Value = 5995.0002
$ws.Range("M137").Value = -2856.75
$ws.Range("N137").Value = -11095.0002

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1782.2572
$ws.Range("J138").Value = 2095.3333
$ws.Range("L138").Value = 6285.999899999999
$ws.Range("N138").Value = -16565.9999

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5334.8696
$ws.Range("I61").Value = 4019.5881
$ws.Range("J61").Value = 9061.5
$ws.Range("K61").Value = 4019.5881
$ws.Range("L61").Value = 9061.5
$ws.Range("M61").Value = -3807.5881
$ws.Range("N61").Value = -9485.5

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 9250728
$ws.Range("I110").Value = 9250728
$ws.Range("K110").Value = 9250728
$ws.Range("M110").Value = -9248683

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5334.8696
$ws.Range("I136").Value = 4019.5881
$ws.Range("J136").Value = 9061.5
$ws.Range("K136").Value = 12058.7643
$ws.Range("L136").Value = 27184.5
$ws.Range("M136").Value = -9508.764299999999
$ws.Range("N136").Value = -32284.5

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1660.8334
$ws.Range("I94").Value = 1693
$ws.Range("K94").Value = 1693
$ws.Range("M94").Value = -1242

# BSM row 103
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 28499.666
$ws.Range("J103").Value = 28499.666
$ws.Range("L103").Value = 28499.666
$ws.Range("N103").Value = -30843.666

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1542.8
$ws.Range("J134").Value = 1949.5
$ws.Range("L134").Value = 5848.5
$ws.Range("N134").Value = -10918.5

# BSM row 140
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 80426.664
$ws.Range("J140").Value = 80426.664
$ws.Range("L140").Value = 80426.664
$ws.Range("N140").Value = -90786.664

# CRP row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 95.411766
$ws.Range("I7").Value = 51.1
$ws.Range("J7").Value = 158.71428
$ws.Range("K7").Value = 51.1
$ws.Range("L7").Value = 158.71428
$ws.Range("M7").Value = 61.9
$ws.Range("N7").Value = -384.71428

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1813.5834
$ws.Range("I31").Value = 1531.3334
$ws.Range("J31").Value = 2660.3333
$ws.Range("K31").Value = 1531.3334
$ws.Range("L31").Value = 2660.3333
$ws.Range("M31").Value = -1236.3334
$ws.Range("N31").Value = -3250.3333

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1813.5834
$ws.Range("I34").Value = 1531.3334
$ws.Range("J34").Value = 2660.3333
$ws.Range("K34").Value = 1531.3334
$ws.Range("L34").Value = 2660.3333
$ws.Range("M34").Value = -1329.3334
$ws.Range("N34").Value = -3064.3333

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1137.75
$ws.Range("I107").Value = 881.5
$ws.Range("K107").Value = 881.5
$ws.Range("M107").Value = 1038.5

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1467.2258
$ws.Range("I132").Value = 1399.4667
$ws.Range("K132").Value = 4198.4001
$ws.Range("M132").Value = -1668.4001

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1004
$ws.Range("I134").Value = 1028.5834
$ws.Range("K134").Value = 3085.7502
$ws.Range("M134").Value = -550.7501999999999

# CUL row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 9167159
$ws.Range("I4").Value = 10000083
$ws.Range("K4").Value = 30000249
$ws.Range("M4").Value = -30000137

# CUL row 24
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2191.4
$ws.Range("I122").Value = 2191.4
$ws.Range("K122").Value = 6574.200000000001
$ws.Range("M122").Value = -4124.200000000001

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3723.25
$ws.Range("I40").Value = 2447
$ws.Range("J40").Value = 4999.5
$ws.Range("K40").Value = 2447
$ws.Range("L40").Value = 4999.5
$ws.Range("M40").Value = -2311
$ws.Range("N40").Value = -5271.5

# LTW row 62
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

# LTW row 65
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

# WVR row 19
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()

# WVR row 49
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

# WVR row 52
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 39999
$ws.Range("J52").Value = 39999
$ws.Range("L52").Value = 39999
$ws.Range("N52").Value = -40451

# WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 13778
$ws.Range("J62").Value = 11166.667
$ws.Range("L62").Value = 11166.667
$ws.Range("N62").Value = -12414.667

# WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 13778
$ws.Range("J65").Value = 11166.667
$ws.Range("L65").Value = 55833.335
$ws.Range("N65").Value = -62073.335

# WVR row 80
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1332.3334
$ws.Range("I81").Value = 1398.5
$ws.Range("J81").Value = 1200
$ws.Range("K81").Value = 2797
$ws.Range("L81").Value = 2400
$ws.Range("M81").Value = -1736
$ws.Range("N81").Value = -4522

# WVR row 83
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1332.3334
$ws.Range("I84").Value = 1398.5
$ws.Range("J84").Value = 1200
$ws.Range("K84").Value = 13985
$ws.Range("L84").Value = 12000
$ws.Range("M84").Value = -8681
$ws.Range("N84").Value = -22608

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2585
$ws.Range("I132").Value = 2695
$ws.Range("J132").Value = 385
$ws.Range("K132").Value = 8085
$ws.Range("L132").Value = 1155
$ws.Range("M132").Value = -5555
$ws.Range("N132").Value = -6215
